$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.860.36"
$ws.Range("E2").Value = "  -4.20%  "
$ws.Range("D3").Value = "3.299.23"
$ws.Range("E3").Value = "  -5.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.36"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.23"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -4.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").Value = "3.290.57"
$ws.Range("E9").Value = "  -4.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.188"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -6.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.589"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.78"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -7.17%  "
$ws.Range("E13").Value = "  -5.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "636.01"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.58"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -5.11%  "
$ws.Range("D16").Value = "3.823.67"
$ws.Range("E16").Value = "  -5.82%  "
$ws.Range("D17").Value = "65.810.16"
$ws.Range("E17").Value = "  -4.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.90"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").Value = "3.297.55"
$ws.Range("E20").Value = "  -5.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.46"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -6.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.906"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -3.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.72"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "107.10"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +7.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.05"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -5.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.99"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -7.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.00"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -5.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.52"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -4.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.74"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -5.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.48"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -5.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.00"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.38"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.11"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "554.44"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +8.75%  "
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.37"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -5.38%  "
$ws.Range("D38").Value = "3.710.50"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.52"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -5.56%  "
$ws.Range("E42").Value = "  -9.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.40"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +34.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.127"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -3.45%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.343"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -5.69%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.21"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -5.28%  "
$ws.Range("E47").Value = "  -5.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.23"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("E49").Value = "  -5.98%  "
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -0.35%  "
